$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 322; this shifts the existing rows 322-349
# down to 323-350 (matching the rest of the diff, which is simply the
# old row N now living at row N+1).
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with the new record.
$ws.Cells.Item(322, 1).Value = 6
$ws.Cells.Item(322, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(322, 3).Value = "Metropolitana"
$ws.Cells.Item(322, 4).Value = 45132
$ws.Cells.Item(322, 5).Value = 13
$ws.Cells.Item(322, 6).Value = 100112029
$ws.Cells.Item(322, 7).Value = "Orégano"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 30
$ws.Cells.Item(322, 11).Value = 20000
$ws.Cells.Item(322, 12).Value = 20000
$ws.Cells.Item(322, 13).Value = 20000
$ws.Cells.Item(322, 14).Value = "$/docena de atados"
$ws.Cells.Item(322, 15).Value = "Región Metropolitana"
$ws.Cells.Item(322, 16).Value = 6667
$ws.Cells.Item(322, 17).Value = 3
$ws.Cells.Item(322, 18).Value = "Hortaliza"
